$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three existing hyperlinks on B1:B3 stay where they are. The hyperlink
# formerly on B4 (Link Noto Emoji 32x32) is logically moved down to B6, and
# two brand-new rows (Link Cerulean 32x32, Link Webuosities 32x32) are
# inserted at B4/B5. This runtime cannot remove or retarget a single
# pre-existing hyperlink cleanly, so clear the whole collection first and
# rebuild every hyperlink (unchanged ones included) afterwards.
$ws.Hyperlinks.Delete()

# Row 4: Link Cerulean 32x32 (new)
$ws.Range("A4").Value = "Link Cerulean 32x32"
$ws.Range("B4").Value = "https://iconarchive.com/show/cerulean-icons-by-iconleak/link-icon.html"

# Row 5: Link Webuosities 32x32 (new). Set the URL cell before the label
# cell so the shared-string table ends up in the same order as the source
# edit.
$ws.Range("B5").Value = "https://iconarchive.com/show/webuosities-icons-by-etherbrian/links-icon.html"
$ws.Range("A5").Value = "Link Webuosities 32x32"

# Row 6: Link Noto Emoji 32x32 (moved down from the old row 4); row 7 stays
# blank.
$ws.Range("A6").Value = "Link Noto Emoji 32x32"
$ws.Range("B6").Value = "https://iconarchive.com/show/noto-emoji-travel-places-icons-by-google/42453-globe-showing-Asia-Australia-icon.html"

# Row 8: explanatory note, plain text only (no hyperlink, default style)
$ws.Range("B8").Value = "These and others have been moved to View\Properties to try to fix Visual Studio resx editor 'file not found' errors"

# Re-create all hyperlinks. Process the moved Noto Emoji link before the two
# brand-new ones so relationship ids come out in the same order as the
# source edit (rId4 keeps pointing at the Noto Emoji target, now from B6).
$ws.Hyperlinks.Add($ws.Range("B1"), "https://www.iconarchive.com/show/flatastic-1-icons-by-custom-icon-design/copy-icon.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.iconarchive.com/show/flatastic-1-icons-by-custom-icon-design/cut-icon.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://www.iconarchive.com/show/flatastic-8-icons-by-custom-icon-design/Paste-icon.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "https://iconarchive.com/show/noto-emoji-travel-places-icons-by-google/42453-globe-showing-Asia-Australia-icon.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://iconarchive.com/show/cerulean-icons-by-iconleak/link-icon.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://iconarchive.com/show/webuosities-icons-by-etherbrian/links-icon.html") | Out-Null

# Hyperlinks.Add can leave a stray unused style behind; force the canonical
# "Hyperlink" cell style back onto every hyperlinked cell.
$ws.Range("B1:B6").Style = "Hyperlink"

# Update the selected cell to match the post-edit state
$ws.Range("B9").Select() | Out-Null
